$d = $word.ActiveDocument

# Locate the paragraph that currently reads:
#   "Initial work involved measuring fatigue life for as-received Ti-6242 alloys. "
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Initial work involved measuring fatigue life for as-received Ti-6242 alloys.*") {
        $target = $cand
        break
    }
}

$paraStart = $target.Range.Start

$firstPhrase  = "Initial work involved measuring fatigue life for "
$secondPhrase = "as-received Ti-6242 alloys"
$thirdPhrase  = " at different fractions of the yield strength"
$fourthPhrase = ". "

$firstStart  = $paraStart
$firstEnd    = $firstStart + $firstPhrase.Length
$secondEnd   = $firstEnd + $secondPhrase.Length

# Replace the old trailing ". " with the new continuation plus the restored ". "
$tailRange = $d.Range($secondEnd, $secondEnd + 2)
$tailRange.Text = $thirdPhrase + $fourthPhrase

$thirdEnd  = $secondEnd + $thirdPhrase.Length
$fourthEnd = $thirdEnd + $fourthPhrase.Length

# Bump every run in the paragraph up to 14pt (sz 28), one segment at a time so
# each original/inserted run keeps its own run-properties block.
$d.Range($firstStart, $firstEnd).Font.Size = 14
$d.Range($firstEnd, $secondEnd).Font.Size = 14
$d.Range($secondEnd, $thirdEnd).Font.Size = 14
$d.Range($thirdEnd, $fourthEnd).Font.Size = 14

# Also bump the paragraph mark itself (w:pPr/w:rPr/w:sz) to match.
$target.Range.Font.Size = 14
